$p = $ppt.ActivePresentation

# 1. Remove the branding/logo pictures from every slide layout (EU funding
#    logos + JetBrains logo) - they were removed from slideLayout1-7.xml.
$design = $p.Designs.Item(1)
$master = $design.SlideMaster
for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $cl = $master.CustomLayouts.Item($i)
    for ($j = $cl.Shapes.Count; $j -ge 1; $j--) {
        $sh = $cl.Shapes.Item($j)
        if ($sh.Type -eq 13) {
            $sh.Delete()
        }
    }
}

# 2. Delete slide 3 ("Termíny lekcí") and slide 2 ("Kurz Java 1" / wifi info)
#    — delete from the back so indices of the remaining slides stay valid.
$p.Slides.Item(3).Delete()
$p.Slides.Item(2).Delete()

# 3. On the remaining "Organizační pokyny" slide (now slide 2), update the
#    course-year URL from 2018 to 2019.
$s = $p.Slides.Item(2)
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        $full = $tr.Text
        $idx = $full.IndexOf("2018-podzim")
        if ($idx -ge 0) {
            $sub = $tr.Characters($idx + 1, 4)
            $sub.Text = "2019"
        }
    }
}
